$d = $word.ActiveDocument

# 1) Split the "Introduce data analytics " paragraph: insert a new paragraph
#    right after it that will hold the new bullet-list item.
$d.Content.Find.Execute(
    "Introduce data analytics ", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Introduce data analytics `rData Analysis Process: ask, prepare, process, analyze, share ",
    2
) | Out-Null

# The newly-created paragraph ("Data Analysis Process: ...") is now the 4th
# paragraph in the document.
$newPara = $d.Paragraphs(4)

# 2) Turn it into a numbered/bulleted list item (this mints word/numbering.xml
#    and wires up numId=1 / abstractNumId=0 on the paragraph).
$newPara.Range.ListFormat.ApplyBulletDefault()

# 3) Apply the "List Paragraph" style (this mints the ListParagraph style in
#    styles.xml with its default shape).
$newPara.Style = "List Paragraph"

# 4) Fix up the minted style so it matches Word's canonical "List Paragraph"
#    quick style: uiPriority 34, 0.5" left indent, contextual spacing.
$listStyle = $d.Styles("List Paragraph")
$listStyle.Priority = 34
$listStyle.ParagraphFormat.LeftIndent = 36
$listStyle.NoSpaceBetweenParagraphsOfSameStyle = $true

# 5) Re-assert the paragraph's exact OOXML (text + pStyle + numPr) so it ends
#    up as a clean <w:p> with no stray rsid bookkeeping attributes.
$xml = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="1"/>
    </w:numPr>
  </w:pPr>
  <w:r>
    <w:t xml:space="preserve">Data Analysis Process: ask, prepare, process, analyze, share </w:t>
  </w:r>
</w:p>
"@
$newPara.Range.InsertXML($xml)
